# Auto-generated Excel COM-interop script applying Mateus_Profits market-data refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H,I,J,K,L,M,N) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 783.87756
$ws.Range("I17").Value = 242
$ws.Range("K17").Value = 726
$ws.Range("M17").Value = -558
$ws.Range("H57").Value = 62999.9
$ws.Range("J57").Value = 29999.857
$ws.Range("L57").Value = 89999.571
$ws.Range("N57").Value = -90997.571
$ws.Range("H87").Value = 58800
$ws.Range("J87").Value = 58800
$ws.Range("L87").Value = 58800
$ws.Range("N87").Value = -61296
$ws.Range("H90").Value = 58800
$ws.Range("J90").Value = 58800
$ws.Range("L90").Value = 176400
$ws.Range("N90").Value = -188880
$ws.Range("H92").Value = 194.66667
$ws.Range("I92").Value = 194.66667
$ws.Range("K92").Value = 194.66667
$ws.Range("M92").Value = 1053.33333
$ws.Range("H98").Value = 2674.3845
$ws.Range("J98").Value = 1549.6666
$ws.Range("L98").Value = 1549.6666
$ws.Range("N98").Value = -4545.6666
$ws.Range("H122").Value = 2674.3845
$ws.Range("J122").Value = 1549.6666
$ws.Range("L122").Value = 4648.9998
$ws.Range("N122").Value = -9548.9998
$ws.Range("H129").Value = 1705.5
$ws.Range("I129").Value = 1447.25
$ws.Range("J129").Value = 2222
$ws.Range("K129").Value = 4341.75
$ws.Range("L129").Value = 6666
$ws.Range("M129").Value = 658.25
$ws.Range("N129").Value = -16666
$ws.Range("H138").Value = 40002430
$ws.Range("J138").Value = 66668684
$ws.Range("L138").Value = 200006052
$ws.Range("N138").Value = -200016332

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 7473.9
$ws.Range("J46").Value = 7854.875
$ws.Range("L46").Value = 7854.875
$ws.Range("N46").Value = -8492.875
$ws.Range("H55").Value = 25137.5
$ws.Range("J55").Value = 35300
$ws.Range("L55").Value = 35300
$ws.Range("N55").Value = -35930
$ws.Range("H61").Value = 13895737
$ws.Range("I61").Value = 16134684
$ws.Range("J61").Value = 14266.6
$ws.Range("K61").Value = 16134684
$ws.Range("L61").Value = 14266.6
$ws.Range("M61").Value = -16134472
$ws.Range("N61").Value = -14690.6
$ws.Range("H136").Value = 13895737
$ws.Range("I136").Value = 16134684
$ws.Range("J136").Value = 14266.6
$ws.Range("K136").Value = 48404052
$ws.Range("L136").Value = 42799.8
$ws.Range("M136").Value = -48401502
$ws.Range("N136").Value = -47899.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 31306
$ws.Range("J35").Value = 31306
$ws.Range("L35").Value = 31306
$ws.Range("N35").Value = -31926
$ws.Range("H82").Value = 43000
$ws.Range("J82").Value = 43000
$ws.Range("L82").Value = 43000
$ws.Range("N82").Value = -43766
$ws.Range("H85").Value = 43000
$ws.Range("J85").Value = 43000
$ws.Range("L85").Value = 43000
$ws.Range("N85").Value = -45652
$ws.Range("H99").Value = 4371.5186
$ws.Range("I99").Value = 2968.4443
$ws.Range("J99").Value = 7177.6665
$ws.Range("K99").Value = 2968.4443
$ws.Range("L99").Value = 7177.6665
$ws.Range("M99").Value = -1470.4443
$ws.Range("N99").Value = -10173.6665
$ws.Range("H105").Value = 3051.4
$ws.Range("I105").Value = 2320.5
$ws.Range("J105").Value = 5975
$ws.Range("K105").Value = 2320.5
$ws.Range("L105").Value = 5975
$ws.Range("M105").Value = -573.5
$ws.Range("N105").Value = -9469

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 110
$ws.Range("J7").Value = 250
$ws.Range("L7").Value = 250
$ws.Range("N7").Value = -476
$ws.Range("H41").Value = 22649.5
$ws.Range("J41").Value = 21100
$ws.Range("L41").Value = 21100
$ws.Range("N41").Value = -21956
$ws.Range("H60").Value = 24400
$ws.Range("J60").Value = 24400
$ws.Range("L60").Value = 24400
$ws.Range("N60").Value = -25422
$ws.Range("H86").Value = 11567.125
$ws.Range("I86").Value = 11720.667
$ws.Range("K86").Value = 11720.667
$ws.Range("M86").Value = -10597.667
$ws.Range("H89").Value = 11567.125
$ws.Range("I89").Value = 11720.667
$ws.Range("K89").Value = 58603.335
$ws.Range("M89").Value = -52987.335
$ws.Range("H132").Value = 3290.375
$ws.Range("I132").Value = 3290.375
$ws.Range("K132").Value = 9871.125
$ws.Range("M132").Value = -7341.125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1499.75
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 6000
$ws.Range("N22").Value = -6338
$ws.Range("H27").Value = 1499.75
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 6000
$ws.Range("N27").Value = -6204
$ws.Range("H32").Value = 474.5
$ws.Range("I32").Value = 399.33334
$ws.Range("J32").Value = 700
$ws.Range("K32").Value = 1198.00002
$ws.Range("L32").Value = 2100
$ws.Range("M32").Value = -915.0000199999999
$ws.Range("N32").Value = -2666
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").Value = $null

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 19333.334
$ws.Range("I18").Value = 19000
$ws.Range("K18").Value = 19000
$ws.Range("M18").Value = -18707
$ws.Range("H80").Value = 3361.8235
$ws.Range("I80").Value = 3430.0667
$ws.Range("K80").Value = 3430.0667
$ws.Range("M80").Value = -2432.0667
$ws.Range("H83").Value = 3361.8235
$ws.Range("I83").Value = 3430.0667
$ws.Range("K83").Value = 17150.3335
$ws.Range("M83").Value = -12158.3335
$ws.Range("H107").Value = 688.82355
$ws.Range("I107").Value = 888.55554
$ws.Range("K107").Value = 888.55554
$ws.Range("M107").Value = 1031.44446

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = $null
$ws.Range("H36").Value = 99999
$ws.Range("J36").Value = 99999
$ws.Range("L36").Value = 99999
$ws.Range("N36").Value = -101123
$ws.Range("H100").Value = 5416.5
$ws.Range("H132").Value = 5194.614
$ws.Range("I132").Value = 5250.1113
$ws.Range("K132").Value = 15750.3339
$ws.Range("M132").Value = -13220.3339

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 730.94446
$ws.Range("I107").Value = 577.26666
$ws.Range("K107").Value = 1731.79998
$ws.Range("M107").Value = 188.20002
$ws.Range("H113").Value = 1226
$ws.Range("J113").Value = 848
$ws.Range("L113").Value = 2544
$ws.Range("N113").Value = -6884
$ws.Range("H126").Value = 2161.077
$ws.Range("I126").Value = 2161.077
$ws.Range("K126").Value = 6483.231000000001
$ws.Range("M126").Value = -4013.231000000001
